$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 12057983.8
$ws.Range("P2").Value = 201.9002445217
$ws.Range("Q2").Value = 342123296.47
$ws.Range("R2").Value = 5728.5511707062
$ws.Range("S2").Value = 78011938.23999999
$ws.Range("T2").Value = 1306.2407171474
$ws.Range("U2").Value = -18430283.89
$ws.Range("V2").Value = -308.5987579445
$ws.Range("W2").Value = 3545588.43
$ws.Range("X2").Value = 59.3677336828
$ws.Range("Y2").Value = 12977896.32
$ws.Range("Z2").Value = 217.3033637999
$ws.Range("AA2").Value = 12344548.33
$ws.Range("AB2").Value = 206.6985134228
$ws.Range("AC2").Value = 5972248.24
$ws.Range("AD2").Value = -40.4451962938
